$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$cell = $ws.Range("D2")
$cell.NumberFormat = "@"
$cell.Value = '63.801.49'
$cell.Style = "Normal"
$ws.Range("E2").Value = '  +1.33%  '

# Row 3
$cell = $ws.Range("D3")
$cell.NumberFormat = "@"
$cell.Value = '2.618.04'
$cell.Style = "Normal"
$ws.Range("E3").Value = '  +0.84%  '

# Row 4
$ws.Range("E4").Value = '  -0.14%  '

# Row 5
$cell = $ws.Range("D5")
$cell.NumberFormat = "@"
$cell.Value = '596.98'
$cell.Style = "Normal"
$ws.Range("E5").Value = '  +2.43%  '

# Row 6
$cell = $ws.Range("D6")
$cell.NumberFormat = "@"
$cell.Value = '146.44'
$cell.Style = "Normal"
$ws.Range("E6").Value = '  -0.69%  '

# Row 7
$cell = $ws.Range("D7")
$cell.NumberFormat = "@"
$cell.Value = '1.00'
$cell.Style = "Normal"
$ws.Range("E7").Value = '  -0.09%  '

# Row 8
$cell = $ws.Range("D8")
$cell.NumberFormat = "@"
$cell.Value = '0.590'
$cell.Style = "Normal"
$ws.Range("E8").Value = '  -1.73%  '

# Row 9
$cell = $ws.Range("D9")
$cell.NumberFormat = "@"
$cell.Value = '0.108'
$cell.Style = "Normal"
$ws.Range("E9").Value = '  -1.29%  '

# Row 10
$cell = $ws.Range("D10")
$cell.NumberFormat = "@"
$cell.Value = '5.65'
$cell.Style = "Normal"
$ws.Range("E10").Value = '  -0.09%  '

# Row 11
$cell = $ws.Range("D11")
$cell.NumberFormat = "@"
$cell.Value = '0.151'
$cell.Style = "Normal"
$ws.Range("E11").Value = '  -0.33%  '

# Row 12
$cell = $ws.Range("D12")
$cell.NumberFormat = "@"
$cell.Value = '0.354'
$cell.Style = "Normal"
$ws.Range("E12").Value = '  -0.85%  '

# Row 13
$cell = $ws.Range("D13")
$cell.NumberFormat = "@"
$cell.Value = '27.40'
$cell.Style = "Normal"
$ws.Range("E13").Value = '  +0.25%  '

# Row 14
$cell = $ws.Range("D14")
$cell.NumberFormat = "@"
$cell.Value = '3.086.11'
$cell.Style = "Normal"
$ws.Range("E14").Value = '  +0.79%  '

# Row 15
$cell = $ws.Range("D15")
$cell.NumberFormat = "@"
$cell.Value = '63.678.88'
$cell.Style = "Normal"
$ws.Range("E15").Value = '  +1.37%  '

# Row 16
$ws.Range("E16").Value = '  -1.84%  '

# Row 17
$cell = $ws.Range("D17")
$cell.NumberFormat = "@"
$cell.Value = '2.592.47'
$cell.Style = "Normal"
$ws.Range("E17").Value = '  -3.76%  '

# Row 18
$cell = $ws.Range("D18")
$cell.NumberFormat = "@"
$cell.Value = '11.22'
$cell.Style = "Normal"
$ws.Range("E18").Value = '  -1.45%  '

# Row 19
$cell = $ws.Range("D19")
$cell.NumberFormat = "@"
$cell.Value = '342.26'
$cell.Style = "Normal"
$ws.Range("E19").Value = '  -0.29%  '

# Row 20
$cell = $ws.Range("D20")
$cell.NumberFormat = "@"
$cell.Value = '4.37'
$cell.Style = "Normal"
$ws.Range("E20").Value = '  -1.10%  '

# Row 21
$cell = $ws.Range("D21")
$cell.NumberFormat = "@"
$cell.Value = '6.74'
$cell.Style = "Normal"
$ws.Range("E21").Value = '  -0.76%  '

# Row 22
$ws.Range("E22").Value = '  +0.76%  '

# Row 23
$cell = $ws.Range("D23")
$cell.NumberFormat = "@"
$cell.Value = '68.93'
$cell.Style = "Normal"
$ws.Range("E23").Value = '  +2.91%  '

# Row 24
$ws.Range("E24").Value = '  +7.27%  '

# Row 25
$cell = $ws.Range("D25")
$cell.NumberFormat = "@"
$cell.Value = '1.63'
$cell.Style = "Normal"
$ws.Range("E25").Value = '  +1.64%  '

# Row 26
$ws.Range("E26").Value = '  -2.40%  '

# Row 27
$ws.Range("B27").Value = 'Aptos'
$ws.Range("C27").Value = 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt'
$cell = $ws.Range("D27")
$cell.NumberFormat = "@"
$cell.Value = '8.00'
$cell.Style = "Normal"
$ws.Range("E27").Value = '  -0.02%  '

# Row 28
$ws.Range("B28").Value = 'Binance-PegBSC-USD'
$ws.Range("C28").Value = 'https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd'
$cell = $ws.Range("D28")
$cell.NumberFormat = "@"
$cell.Value = '1.00'
$cell.Style = "Normal"
$ws.Range("E28").Value = '  +1.71%  '

# Row 29
$cell = $ws.Range("D29")
$cell.NumberFormat = "@"
$cell.Value = '8.37'
$cell.Style = "Normal"
$ws.Range("E29").Value = '  -1.03%  '

# Row 30
$cell = $ws.Range("D30")
$cell.NumberFormat = "@"
$cell.Value = '2.03'
$cell.Style = "Normal"
$ws.Range("E30").Value = '  +5.25%  '

# Row 31
$cell = $ws.Range("D31")
$cell.NumberFormat = "@"
$cell.Value = '498.15'
$cell.Style = "Normal"
$ws.Range("E31").Value = '  +8.25%  '

# Row 32
$cell = $ws.Range("D32")
$cell.NumberFormat = "@"
$cell.Value = '1.73'
$cell.Style = "Normal"
$ws.Range("E32").Value = '  +6.44%  '

# Row 33
$cell = $ws.Range("D33")
$cell.NumberFormat = "@"
$cell.Value = '0.0₃0816'
$cell.Style = "Normal"
$ws.Range("E33").Value = '  -1.12%  '

# Row 34
$cell = $ws.Range("D34")
$cell.NumberFormat = "@"
$cell.Value = '174.61'
$cell.Style = "Normal"
$ws.Range("E34").Value = '  -1.11%  '

# Row 35
$ws.Range("E35").Value = '  -0.13%  '

# Row 36
$cell = $ws.Range("D36")
$cell.NumberFormat = "@"
$cell.Value = '0.403'
$cell.Style = "Normal"
$ws.Range("E36").Value = '  -1.51%  '

# Row 37
$cell = $ws.Range("D37")
$cell.NumberFormat = "@"
$cell.Value = '19.06'
$cell.Style = "Normal"
$ws.Range("E37").Value = '  -1.15%  '

# Row 38
$ws.Range("E38").Value = '  -0.85%  '

# Row 39
$ws.Range("B39").Value = 'Stacks'
$ws.Range("C39").Value = 'https://coinranking.com/coin/mMPrMcB7+stacks-stx'
$cell = $ws.Range("D39")
$cell.NumberFormat = "@"
$cell.Value = '1.74'
$cell.Style = "Normal"
$ws.Range("E39").Value = '  +1.43%  '

# Row 40
$ws.Range("B40").Value = 'USDe'
$ws.Range("C40").Value = 'https://coinranking.com/coin/exbfr2U-0+usde-usde'
$cell = $ws.Range("D40")
$cell.NumberFormat = "@"
$cell.Value = '0.999'
$cell.Style = "Normal"
$ws.Range("E40").Value = '  +0.04%  '

# Row 41
$cell = $ws.Range("D41")
$cell.NumberFormat = "@"
$cell.Value = '166.20'
$cell.Style = "Normal"
$ws.Range("E41").Value = '  +4.74%  '

# Row 42
$cell = $ws.Range("D42")
$cell.NumberFormat = "@"
$cell.Value = '40.16'
$cell.Style = "Normal"
$ws.Range("E42").Value = '  +1.56%  '

# Row 43
$cell = $ws.Range("D43")
$cell.NumberFormat = "@"
$cell.Value = '3.77'
$cell.Style = "Normal"
$ws.Range("E43").Value = '  -0.95%  '

# Row 44
$cell = $ws.Range("D44")
$cell.NumberFormat = "@"
$cell.Value = '21.83'
$cell.Style = "Normal"
$ws.Range("E44").Value = '  +5.33%  '

# Row 45
$cell = $ws.Range("D45")
$cell.NumberFormat = "@"
$cell.Value = '0.628'
$cell.Style = "Normal"
$ws.Range("E45").Value = '  -1.27%  '

# Row 46
$cell = $ws.Range("D46")
$cell.NumberFormat = "@"
$cell.Value = '0.0544'
$cell.Style = "Normal"
$ws.Range("E46").Value = '  -0.39%  '

# Row 47
$ws.Range("B47").Value = 'VeChain'
$ws.Range("C47").Value = 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet'
$cell = $ws.Range("D47")
$cell.NumberFormat = "@"
$cell.Value = '0.0240'
$cell.Style = "Normal"
$ws.Range("E47").Value = '  +0.77%  '

# Row 48
$ws.Range("B48").Value = 'Stellar'
$ws.Range("C48").Value = 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm'
$cell = $ws.Range("D48")
$cell.NumberFormat = "@"
$cell.Value = '0.0960'
$cell.Style = "Normal"
$ws.Range("E48").Value = '  -1.39%  '

# Row 49
$cell = $ws.Range("D49")
$cell.NumberFormat = "@"
$cell.Value = '18.61'
$cell.Style = "Normal"
$ws.Range("E49").Value = '  -0.25%  '

# Row 50
$cell = $ws.Range("D50")
$cell.NumberFormat = "@"
$cell.Value = '1.75'
$cell.Style = "Normal"
$ws.Range("E50").Value = '  +0.43%  '

# Row 51
$cell = $ws.Range("D51")
$cell.NumberFormat = "@"
$cell.Value = '11.36'
$cell.Style = "Normal"
$ws.Range("E51").Value = '  -0.44%  '
